$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (item id 27772)
$ws.Range("H28").Value = 1189.5
$ws.Range("I28").Value = 685.1539
$ws.Range("K28").Value = 685.1539
$ws.Range("M28").Value = -200.1539

# Row 54 (item id 2174)
$ws.Range("H54").Value = 25000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 25000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 25000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -25972

# Row 98 (item id 36237)
$ws.Range("H98").Value = 17617.5
$ws.Range("I98").Value = 1051.25
$ws.Range("J98").Value = 50750
$ws.Range("K98").Value = 1051.25
$ws.Range("L98").Value = 50750
$ws.Range("M98").Value = 446.75
$ws.Range("N98").Value = -53746

# Row 101 (item id 19884)
$ws.Range("H101").Value = 1229.9
$ws.Range("I101").Value = 759.7143
$ws.Range("J101").Value = 2327
$ws.Range("K101").Value = 2279.1429
$ws.Range("L101").Value = 6981
$ws.Range("M101").Value = -657.1428999999998
$ws.Range("N101").Value = -10225

# Row 118 (item id 27958)
$ws.Range("H118").Value = 900.5
$ws.Range("I118").Value = 253.5
$ws.Range("K118").Value = 760.5
$ws.Range("M118").Value = 896.5

# Row 122 (item id 36237)
$ws.Range("H122").Value = 17617.5
$ws.Range("I122").Value = 1051.25
$ws.Range("J122").Value = 50750
$ws.Range("K122").Value = 3153.75
$ws.Range("L122").Value = 152250
$ws.Range("M122").Value = -703.75
$ws.Range("N122").Value = -157150

# Row 132 (item id 44049)
$ws.Range("H132").Value = 1649.2333
$ws.Range("I132").Value = 1479
$ws.Range("K132").Value = 4437
$ws.Range("M132").Value = -1907

# Row 138 (item id 44169)
$ws.Range("H138").Value = 3755.261
$ws.Range("I138").Value = 3281.077
$ws.Range("J138").Value = 3942.0605
$ws.Range("K138").Value = 9843.231
$ws.Range("L138").Value = 11826.1815
$ws.Range("M138").Value = -4703.231
$ws.Range("N138").Value = -22106.1815

$ws = $wb.Worksheets.Item("ARM")
# Row 97 (item id 19941)
$ws.Range("H97").Value = 178.5
$ws.Range("I97").Value = 175.28572
$ws.Range("K97").Value = 175.28572
$ws.Range("M97").Value = 320.71428

# Row 102 (item id 19945)
$ws.Range("H102").Value = 1704
$ws.Range("I102").Value = 1784.5333
$ws.Range("K102").Value = 1784.5333
$ws.Range("M102").Value = -162.5333000000001

# Row 109 (item id 25646)
$ws.Range("H109").Value = 78947
$ws.Range("J109").Value = 78947
$ws.Range("L109").Value = 78947
$ws.Range("N109").Value = -81721

# Row 119 (item id 26287)
$ws.Range("H119").Value = 65876.39999999999
$ws.Range("J119").Value = 65876.39999999999
$ws.Range("L119").Value = 65876.39999999999
$ws.Range("N119").Value = -75552.39999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 9 (item id 1648)
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20336

# Row 86 (item id 12526)
$ws.Range("H86").Value = 2967.647
$ws.Range("J86").Value = 6126.75
$ws.Range("L86").Value = 6126.75
$ws.Range("N86").Value = -8372.75

# Row 89 (item id 12526)
$ws.Range("H89").Value = 2967.647
$ws.Range("J89").Value = 6126.75
$ws.Range("L89").Value = 30633.75
$ws.Range("N89").Value = -41865.75

# Row 133 (item id 43209)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (item id 5367)
$ws.Range("H22").Value = 1553.9445
$ws.Range("I22").Value = 399.53845
$ws.Range("J22").Value = 4555.4
$ws.Range("K22").Value = 399.53845
$ws.Range("L22").Value = 4555.4
$ws.Range("M22").Value = -49.53845000000001
$ws.Range("N22").Value = -5255.4

# Row 86 (item id 12584)
$ws.Range("H86").Value = 11133.277
$ws.Range("I86").Value = 8399.272000000001
$ws.Range("K86").Value = 8399.272000000001
$ws.Range("M86").Value = -7276.272000000001

# Row 89 (item id 12584)
$ws.Range("H89").Value = 11133.277
$ws.Range("I89").Value = 8399.272000000001
$ws.Range("K89").Value = 41996.36
$ws.Range("M89").Value = -36380.36

# Row 107 (item id 27689)
$ws.Range("H107").Value = 1083
$ws.Range("I107").Value = 1014.5833
$ws.Range("J107").Value = 1200.2858
$ws.Range("K107").Value = 1014.5833
$ws.Range("L107").Value = 1200.2858
$ws.Range("M107").Value = 905.4167
$ws.Range("N107").Value = -5040.2858

# Row 122 (item id 36196)
$ws.Range("H122").Value = 8315.3125
$ws.Range("I122").Value = 2574.625
$ws.Range("K122").Value = 7723.875
$ws.Range("M122").Value = -5273.875

# Row 132 (item id 44019)
$ws.Range("H132").Value = 3819.9583
$ws.Range("I132").Value = 5730.3
$ws.Range("J132").Value = 2455.4285
$ws.Range("K132").Value = 17190.9
$ws.Range("L132").Value = 7366.2855
$ws.Range("M132").Value = -14660.9
$ws.Range("N132").Value = -12426.2855

# Row 134 (item id 44020)
$ws.Range("H134").Value = 2148.2903
$ws.Range("I134").Value = 1381.409
$ws.Range("J134").Value = 4022.889
$ws.Range("K134").Value = 4144.227000000001
$ws.Range("L134").Value = 12068.667
$ws.Range("M134").Value = -1609.227000000001
$ws.Range("N134").Value = -17138.667

$ws = $wb.Worksheets.Item("CUL")
# Row 129 (item id 36054)
$ws.Range("H129").Value = 5560194
$ws.Range("J129").Value = 6949626.5
$ws.Range("L129").Value = 20848879.5
$ws.Range("N129").Value = -20858879.5

$ws = $wb.Worksheets.Item("GSM")
# Row 24 (item id 4431)
$ws.Range("H24").Value = 22383.8
$ws.Range("I24").Value = 24304
$ws.Range("J24").Value = 19503.5
$ws.Range("K24").Value = 24304
$ws.Range("L24").Value = 19503.5
$ws.Range("M24").Value = -24131
$ws.Range("N24").Value = -19849.5

# Row 70 (item id 14146)
$ws.Range("H70").Value = 15477.523
$ws.Range("I70").Value = 5939.5625
$ws.Range("K70").Value = 5939.5625
$ws.Range("M70").Value = -5669.5625

# Row 73 (item id 14146)
$ws.Range("H73").Value = 15477.523
$ws.Range("I73").Value = 5939.5625
$ws.Range("K73").Value = 5939.5625
$ws.Range("M73").Value = -5003.5625

# Row 102 (item id 36169)
$ws.Range("H102").Value = 2443.2307
$ws.Range("I102").Value = 2443.2307
$ws.Range("K102").Value = 2443.2307
$ws.Range("M102").Value = -821.2307000000001

# Row 118 (item id 26172)
$ws.Range("H118").Value = 23875.25
$ws.Range("J118").Value = 23875.25
$ws.Range("L118").Value = 23875.25
$ws.Range("N118").Value = -27189.25

# Row 122 (item id 36182)
$ws.Range("H122").Value = 5378.067
$ws.Range("I122").Value = 4969.4546
$ws.Range("J122").Value = 6501.75
$ws.Range("K122").Value = 14908.3638
$ws.Range("L122").Value = 19505.25
$ws.Range("M122").Value = -12458.3638
$ws.Range("N122").Value = -24405.25

# Row 132 (item id 44008)
$ws.Range("H132").Value = 5487.9287
$ws.Range("I132").Value = 3264.0476
$ws.Range("K132").Value = 9792.1428
$ws.Range("M132").Value = -7262.1428

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (item id 36249)
$ws.Range("H7").Value = 7909.2104
$ws.Range("I7").Value = 4057.5715
$ws.Range("K7").Value = 4057.5715
$ws.Range("M7").Value = -3945.5715

# Row 16 (item id 5289)
$ws.Range("H16").Value = 2939.5334
$ws.Range("I16").Value = 2542.3572
$ws.Range("J16").Value = 8500
$ws.Range("K16").Value = 2542.3572
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = -2372.3572
$ws.Range("N16").Value = -8840

# Row 20 (item id 4308)
$ws.Range("H20").Value = 26501.25
$ws.Range("I20").Value = 25002.5
$ws.Range("J20").Value = 28000
$ws.Range("K20").Value = 25002.5
$ws.Range("L20").Value = 28000
$ws.Range("M20").Value = -24776.5
$ws.Range("N20").Value = -28452

# Row 22 (item id 5277)
$ws.Range("H22").Value = 4883.25
$ws.Range("I22").Value = 2108.5
$ws.Range("J22").Value = 7658
$ws.Range("K22").Value = 2108.5
$ws.Range("L22").Value = 7658
$ws.Range("M22").Value = -1813.5
$ws.Range("N22").Value = -8248

# Row 27 (item id 5277)
$ws.Range("H27").Value = 4883.25
$ws.Range("I27").Value = 2108.5
$ws.Range("J27").Value = 7658
$ws.Range("K27").Value = 2108.5
$ws.Range("L27").Value = 7658
$ws.Range("M27").Value = -2001.5
$ws.Range("N27").Value = -7872

# Row 40 (item id 36248)
$ws.Range("H40").Value = 8706.214
$ws.Range("I40").Value = 7534.273
$ws.Range("J40").Value = 13003.333
$ws.Range("K40").Value = 7534.273
$ws.Range("L40").Value = 13003.333
$ws.Range("M40").Value = -7398.273
$ws.Range("N40").Value = -13275.333

# Row 46 (item id 5282)
$ws.Range("H46").Value = 6500.3335
$ws.Range("J46").Value = 6500.3335
$ws.Range("L46").Value = 6500.3335
$ws.Range("N46").Value = -6876.3335

# Row 61 (item id 27740)
$ws.Range("H61").Value = 6328.4287
$ws.Range("I61").Value = 1433.3334
$ws.Range("K61").Value = 1433.3334
$ws.Range("M61").Value = -1231.3334

# Row 108 (item id 25655)
$ws.Range("H108").Value = 73520
$ws.Range("J108").Value = 73520
$ws.Range("L108").Value = 73520
$ws.Range("N108").Value = -81200

# Row 113 (item id 27740)
$ws.Range("H113").Value = 6328.4287
$ws.Range("I113").Value = 1433.3334
$ws.Range("K113").Value = 1433.3334
$ws.Range("M113").Value = 736.6666

# Row 122 (item id 36247)
$ws.Range("H122").Value = 13001.667
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 19005
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 57015
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -61915

# Row 126 (item id 36249)
$ws.Range("H126").Value = 7909.2104
$ws.Range("I126").Value = 4057.5715
$ws.Range("K126").Value = 12172.7145
$ws.Range("M126").Value = -9702.7145

# Row 132 (item id 44058)
$ws.Range("H132").Value = 4737.409
$ws.Range("I132").Value = 2204
$ws.Range("J132").Value = 13351
$ws.Range("K132").Value = 6612
$ws.Range("L132").Value = 40053
$ws.Range("M132").Value = -4082
$ws.Range("N132").Value = -45113

$ws = $wb.Worksheets.Item("WVR")
# Row 14 (item id 2658)
$ws.Range("H14").Value = 16666
$ws.Range("I14").Value = 14999
$ws.Range("K14").Value = 14999
$ws.Range("M14").Value = -14831

# Row 15 (item id 2670)
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 122 (item id 36208)
$ws.Range("H122").Value = 6601.6
$ws.Range("I122").Value = 2236.3794
$ws.Range("K122").Value = 6709.138199999999
$ws.Range("M122").Value = -4259.138199999999

# Row 126 (item id 36210)
$ws.Range("H126").Value = 2915.0833
$ws.Range("I126").Value = 1399.4286
$ws.Range("J126").Value = 5037
$ws.Range("K126").Value = 4198.2858
$ws.Range("L126").Value = 15111
$ws.Range("M126").Value = -1728.2858
$ws.Range("N126").Value = -20051

# Row 132 (item id 44029)
$ws.Range("H132").Value = 5357.967
$ws.Range("I132").Value = 4738.3335
$ws.Range("J132").Value = 10934.667
$ws.Range("K132").Value = 14215.0005
$ws.Range("L132").Value = 32804.001
$ws.Range("M132").Value = -11685.0005
$ws.Range("N132").Value = -37864.001
